$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 245.33333
$ws.Range("I5").Value = 94.59999999999999
$ws.Range("K5").Value = 94.59999999999999
$ws.Range("M5").Value = 20.40000000000001
$ws.Range("H17").Value = 1904.6086
$ws.Range("J17").Value = 2300.375
$ws.Range("L17").Value = 6901.125
$ws.Range("N17").Value = -7237.125
$ws.Range("H28").Value = 7376.9
$ws.Range("I28").Value = 1176.6666
$ws.Range("J28").Value = 10034.143
$ws.Range("K28").Value = 1176.6666
$ws.Range("L28").Value = 10034.143
$ws.Range("M28").Value = -691.6666
$ws.Range("N28").Value = -11004.143
$ws.Range("H32").Value = 250
$ws.Range("I32").Value = 250
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 250
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 76
$ws.Range("N32").ClearContents()
$ws.Range("H40").Value = 5787.25
$ws.Range("I40").Value = 4499.375
$ws.Range("J40").Value = 8363
$ws.Range("K40").Value = 4499.375
$ws.Range("L40").Value = 8363
$ws.Range("M40").Value = -4324.375
$ws.Range("N40").Value = -8713
$ws.Range("H62").Value = 11299.5
$ws.Range("J62").Value = 11444
$ws.Range("L62").Value = 11444
$ws.Range("N62").Value = -12692
$ws.Range("H65").Value = 11299.5
$ws.Range("J65").Value = 11444
$ws.Range("L65").Value = 57220
$ws.Range("N65").Value = -63460
$ws.Range("H106").Value = 7389.6
$ws.Range("I106").Value = 7237.25
$ws.Range("K106").Value = 7237.25
$ws.Range("M106").Value = -6606.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2775.8462
$ws.Range("I45").Value = 1908.6
$ws.Range("J45").Value = 5666.6665
$ws.Range("K45").Value = 1908.6
$ws.Range("L45").Value = 5666.6665
$ws.Range("M45").Value = -1531.6
$ws.Range("N45").Value = -6420.6665
$ws.Range("H119").Value = 24698
$ws.Range("J119").Value = 24698
$ws.Range("L119").Value = 24698
$ws.Range("N119").Value = -34374

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 25910.666
$ws.Range("I26").Value = 25910.666
$ws.Range("K26").Value = 25910.666
$ws.Range("M26").Value = -25618.666
$ws.Range("H86").Value = 6965
$ws.Range("I86").Value = 2115.8
$ws.Range("J86").Value = 9995.75
$ws.Range("K86").Value = 2115.8
$ws.Range("L86").Value = 9995.75
$ws.Range("M86").Value = -992.8000000000002
$ws.Range("N86").Value = -12241.75
$ws.Range("H89").Value = 6965
$ws.Range("I89").Value = 2115.8
$ws.Range("J89").Value = 9995.75
$ws.Range("K89").Value = 10579
$ws.Range("L89").Value = 49978.75
$ws.Range("M89").Value = -4963
$ws.Range("N89").Value = -61210.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 4162.25
$ws.Range("I25").Value = 4162.25
$ws.Range("K25").Value = 4162.25
$ws.Range("M25").Value = -3988.25
$ws.Range("H58").Value = 2022.3889
$ws.Range("I58").Value = 1094.2667
$ws.Range("K58").Value = 1094.2667
$ws.Range("M58").Value = -891.2666999999999
$ws.Range("H62").Value = 1999.3334
$ws.Range("I62").Value = 2249
$ws.Range("J62").Value = 1500
$ws.Range("K62").Value = 2249
$ws.Range("L62").Value = 1500
$ws.Range("M62").Value = -1625
$ws.Range("N62").Value = -2748
$ws.Range("H65").Value = 1999.3334
$ws.Range("I65").Value = 2249
$ws.Range("J65").Value = 1500
$ws.Range("K65").Value = 11245
$ws.Range("L65").Value = 7500
$ws.Range("M65").Value = -8125
$ws.Range("N65").Value = -13740
$ws.Range("H122").Value = 915.6
$ws.Range("I122").Value = 915.6
$ws.Range("K122").Value = 2746.8
$ws.Range("M122").Value = -296.8000000000002
$ws.Range("H132").Value = 5076.121
$ws.Range("I132").Value = 4404.294
$ws.Range("J132").Value = 5789.9375
$ws.Range("K132").Value = 13212.882
$ws.Range("L132").Value = 17369.8125
$ws.Range("M132").Value = -10682.882
$ws.Range("N132").Value = -22429.8125
$ws.Range("H136").Value = 2022.3889
$ws.Range("I136").Value = 1094.2667
$ws.Range("K136").Value = 3282.800099999999
$ws.Range("M136").Value = -732.8000999999995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 690.8
$ws.Range("I5").Value = 690.8
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2072.4
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1960.4
$ws.Range("N5").ClearContents()
$ws.Range("H8").Value = 365.75
$ws.Range("I8").Value = 365.75
$ws.Range("K8").Value = 1097.25
$ws.Range("M8").Value = -958.25
$ws.Range("H132").Value = 2226.8125
$ws.Range("J132").Value = 2312.1
$ws.Range("L132").Value = 20808.9
$ws.Range("N132").Value = -25868.9
$ws.Range("H135").Value = 690.8
$ws.Range("I135").Value = 690.8
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 6217.2
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -3682.2
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 23333
$ws.Range("J63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("N63").Value = -21372
$ws.Range("H66").Value = 23333
$ws.Range("J66").Value = 20000
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -66864

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 44444
$ws.Range("I63").Value = 44444
$ws.Range("K63").Value = 44444
$ws.Range("M63").Value = -43695
$ws.Range("H66").Value = 44444
$ws.Range("I66").Value = 44444
$ws.Range("K66").Value = 133332
$ws.Range("M66").Value = -129588
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 1752
$ws.Range("I136").Value = 1500
$ws.Range("M136").Value = -1950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 2950
$ws.Range("I51").Value = 2900
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 2900
$ws.Range("L51").Value = 3000
$ws.Range("M51").Value = -2390
$ws.Range("N51").Value = -4020
$ws.Range("H64").Value = 49999.5
$ws.Range("J64").Value = 49999.5
$ws.Range("L64").Value = 49999.5
$ws.Range("N64").Value = -50495.5
$ws.Range("H67").Value = 49999.5
$ws.Range("J67").Value = 49999.5
$ws.Range("L67").Value = 49999.5
$ws.Range("N67").Value = -51715.5
$ws.Range("H75").Value = 34000
$ws.Range("I75").Value = 34000
$ws.Range("K75").Value = 34000
$ws.Range("M75").Value = -33064
$ws.Range("H78").Value = 34000
$ws.Range("I78").Value = 34000
$ws.Range("K78").Value = 102000
$ws.Range("M78").Value = -97320
$ws.Range("H107").Value = 546.125
$ws.Range("I107").Value = 494.36365
$ws.Range("J107").Value = 660
$ws.Range("K107").Value = 1483.09095
$ws.Range("L107").Value = 1980
$ws.Range("M107").Value = 436.90905
$ws.Range("N107").Value = -5820
$ws.Range("H122").Value = 2422.7646
$ws.Range("I122").Value = 1156.5714
$ws.Range("J122").Value = 3309.1
$ws.Range("K122").Value = 3469.7142
$ws.Range("L122").Value = 9927.299999999999
$ws.Range("M122").Value = -1019.7142
$ws.Range("N122").Value = -14827.3
$ws.Range("H132").Value = 2803.077
$ws.Range("I132").Value = 2153.6365
$ws.Range("K132").Value = 6460.9095
$ws.Range("M132").Value = -3930.9095
$ws.Range("H136").Value = 4626.0527
$ws.Range("I136").Value = 3181.125
$ws.Range("K136").Value = 9543.375
$ws.Range("M136").Value = -6993.375
